# Daily attendance processing - 2025-11-28 19:22:51
# For every row in column G ("Recorded By") whose value begins with the
# token "System", reverse the order of the comma-separated recorder list
# (e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -ne $value -and $value -is [string]) {
        $parts = $value -split ',\s*'
        if ($parts.Length -gt 1 -and $parts[0] -eq 'System') {
            $n = $parts.Length
            $reversed = $parts[($n - 1)..0]
            $cell.Value2 = [string]::Join(', ', $reversed)
        }
    }
}
